# Update NATMI ligand-receptor pair statistics for Tgfb2-Tgfbr2
# (expressing-cell counts changed from 1 to 3 for both ligand and receptor,
# which cascades into the derived expression/specificity/edge-weight columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "E"=3; "G"=2.191602; "H"=6.574806000000001; "I"=0.07674610985252207; "J"=0.07674610985252209; "K"=3; "M"=19.60726733333334; "N"=58.82180200000001; "O"=0.1509859438163708; "P"=0.1509859438163708; "Q"=42.971326302268; "R"=386.7419367204121; "S"=0.01158758383031792; "T"=0.01158758383031792 }
    3 = @{ "E"=3; "G"=2.191602; "H"=6.574806000000001; "I"=0.07674610985252207; "J"=0.07674610985252209; "K"=3; "M"=70.710031; "N"=212.130093; "O"=0.544503249041223; "P"=0.544503249041223; "Q"=154.968245359662; "R"=1394.714208236958; "S"=0.04178850616597289; "T"=0.0417885061659729 }
    4 = @{ "E"=3; "G"=2.191602; "H"=6.574806000000001; "I"=0.07674610985252207; "J"=0.07674610985252209; "K"=3; "M"=29.95517733333334; "N"=89.865532; "O"=0.2306701206736283; "P"=0.2306701206736284; "Q"=65.64982655408801; "R"=590.8484389867921; "S"=0.0177030344209128; "T"=0.01770303442091281 }
    5 = @{ "E"=3; "G"=2.191602; "H"=6.574806000000001; "I"=0.07674610985252207; "J"=0.07674610985252209; "K"=3; "M"=9.589065333333332; "N"=28.767196; "O"=0.07384068646877778; "P"=0.0738406864687778; "Q"=21.015414762664; "R"=189.138732863976; "S"=0.00566698543531846; "T"=0.005666985435318462 }
    6 = @{ "E"=3; "G"=16.59481266666667; "H"=49.78443799999999; "I"=0.5811216251390648; "J"=0.5811216251390647; "K"=3; "M"=19.60726733333334; "N"=58.82180200000001; "O"=0.1509859438163708; "P"=0.1509859438163708; "Q"=325.3789283019196; "R"=2928.410354717276; "S"=0.08774119704372495; "T"=0.08774119704372493 }
    7 = @{ "E"=3; "G"=16.59481266666667; "H"=49.78443799999999; "I"=0.5811216251390648; "J"=0.5811216251390647; "K"=3; "M"=70.710031; "N"=212.130093; "O"=0.544503249041223; "P"=0.544503249041223; "Q"=1173.419718099193; "R"=10560.77746289273; "S"=0.3164226129763365; "T"=0.3164226129763364 }
    8 = @{ "E"=3; "G"=16.59481266666667; "H"=49.78443799999999; "I"=0.5811216251390648; "J"=0.5811216251390647; "K"=3; "M"=29.95517733333334; "N"=89.865532; "O"=0.2306701206736283; "P"=0.2306701206736284; "Q"=497.1005562434462; "R"=4473.905006191016; "S"=0.1340473953968831; "T"=0.1340473953968831 }
    9 = @{ "E"=3; "G"=16.59481266666667; "H"=49.78443799999999; "I"=0.5811216251390648; "J"=0.5811216251390647; "K"=3; "M"=9.589065333333332; "N"=28.767196; "O"=0.07384068646877778; "P"=0.0738406864687778; "Q"=159.1287428550942; "R"=1432.158685695848; "S"=0.0429104197221203; "T"=0.0429104197221203 }
    10 = @{ "E"=3; "G"=9.770107666666666; "H"=29.310323; "I"=0.342132265008413; "J"=0.342132265008413; "K"=3; "M"=19.60726733333334; "N"=58.82180200000001; "O"=0.1509859438163708; "P"=0.1509859438163708; "Q"=191.5651128957829; "R"=1724.086016062046; "S"=0.05165716294232794; "T"=0.05165716294232794 }
    11 = @{ "E"=3; "G"=9.770107666666666; "H"=29.310323; "I"=0.342132265008413; "J"=0.342132265008413; "K"=3; "M"=70.710031; "N"=212.130093; "O"=0.544503249041223; "P"=0.544503249041223; "Q"=690.8446159833377; "R"=6217.601543850038; "S"=0.1862921298989136; "T"=0.1862921298989136 }
    12 = @{ "E"=3; "G"=9.770107666666666; "H"=29.310323; "I"=0.342132265008413; "J"=0.342132265008413; "K"=3; "M"=29.95517733333334; "N"=89.865532; "O"=0.2306701206736283; "P"=0.2306701206736284; "Q"=292.6653077207596; "R"=2633.987769486836; "S"=0.07891969085583242; "T"=0.07891969085583243 }
    13 = @{ "E"=3; "G"=9.770107666666666; "H"=29.310323; "I"=0.342132265008413; "J"=0.342132265008413; "K"=3; "M"=9.589065333333332; "N"=28.767196; "O"=0.07384068646877778; "P"=0.0738406864687778; "Q"=93.68620072936754; "R"=843.1758065643079; "S"=0.02526328131133902; "T"=0.02526328131133902 }
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    foreach ($col in $vals.Keys) {
        $ws.Range("$col$r").Value = $vals[$col]
    }
}
